$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet currently has a text header in row 1 (Lg., Threading, HeadDia., ...)
# followed by data rows. The edit inserts a brand-new row 1 containing the
# numeric column indexes (0..10), pushing the old header down to row 2 and
# all the data rows down by one (row 2 -> row 3, ..., row 61 -> row 62).
# The old header row (now row 2) loses its package-code/price/thread-size
# values in columns H, J, K (left blank) while keeping the rest of its text
# labels, and it also loses the bold/border formatting that row 1 used to
# have - that formatting now belongs solely to the new numeric row 1.

# Insert a new blank row above row 1; everything currently in row 1 (and
# below) shifts down by one row.
$ws.Rows.Item(1).Insert()

# The shifted-down old header (now row 2) carried its original formatting
# (bold font, borders, centered) down with it. Grab that formatting for the
# new row 1 first...
$ws.Range("A2:K2").Copy()
$ws.Range("A1:K1").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$excel.CutCopyMode = 0

# ...then reset row 2 back to the plain/default formatting used by the rest
# of the data rows (e.g. row 3, which is untouched data).
$ws.Range("A3:K3").Copy()
$ws.Range("A2:K2").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$excel.CutCopyMode = 0

# Populate the new row 1 with the numeric column indexes 0..10.
for ($col = 1; $col -le 11; $col++) {
    $ws.Cells.Item(1, $col).Value = $col - 1
}

# Clear the package code / price / thread-size cells on the old header row
# (now row 2), columns H, J, K.
$ws.Cells.Item(2, 8).Value = ""
$ws.Cells.Item(2, 10).Value = ""
$ws.Cells.Item(2, 11).Value = ""
